# Added categories and sheet for ventilators
#
# 1. lsh_sheet_names: new raw sheet-name entry "Öndunarvélar - tímar"
# 2. lsh_unit_categories: two new unit-category rows
#      - Fv-B3 GD Æðaskurðlækninga  -> Göngudeild / outpatient_clinic / home
#      - Sjúkrahótel Landspítala (Ám9) -> Sjúkrahótel / patient_hotel / home
# 3. Selection/active-sheet state updated to match the edited sheets.

$wb = $excel.ActiveWorkbook

# --- lsh_sheet_names: add the new raw sheet name row first so the shared
#     string table gets the new strings in the same order as the source
#     workbook (sheet names table, then unit categories table). ---
$wsNames = $wb.Worksheets.Item("lsh_sheet_names")
$wsNames.Range("A12").Value = "Öndunarvélar - tímar"

# --- lsh_unit_categories: append the two new category rows. ---
$wsUnits = $wb.Worksheets.Item("lsh_unit_categories")
$wsUnits.Range("A21").Value = "Fv-B3 GD Æðaskurðlækninga"
$wsUnits.Range("B21").Value = "Göngudeild"
$wsUnits.Range("C21").Value = "outpatient_clinic"
$wsUnits.Range("D21").Value = "home"
$wsUnits.Range("E21").Value = 1

$wsUnits.Range("A22").Value = "Sjúkrahótel Landspítala (Ám9)"
$wsUnits.Range("B22").Value = "Sjúkrahótel"
$wsUnits.Range("C22").Value = "patient_hotel"
$wsUnits.Range("D22").Value = "home"
$wsUnits.Range("E22").Value = 1

# --- View/selection bookkeeping matching the authored workbook state. ---
$wsNames.Activate()
$wsNames.Range("A12").Select()

$wsUnits.Activate()
$wsUnits.Range("B23").Select()
